# Update the sensor data (columns C-H) for rows 2-21 using the sliding
# window of walkingToRunning samples, and append new rows 22-31 that
# extend the dataset (column A/B pattern continues: timestamp +100ms, same label).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "walkingToRunning"
$ws.Cells.Item(2, 3).Value = 4.936917678169539
$ws.Cells.Item(2, 4).Value = -12.71771009072015
$ws.Cells.Item(2, 5).Value = 3.41959030731864
$ws.Cells.Item(2, 6).Value = 2.903977394104004
$ws.Cells.Item(2, 7).Value = -3.523428678512573
$ws.Cells.Item(2, 8).Value = -6.097240924835205

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "walkingToRunning"
$ws.Cells.Item(3, 3).Value = 11.45307619675343
$ws.Cells.Item(3, 4).Value = -9.954130877619217
$ws.Cells.Item(3, 5).Value = 7.623632503592415
$ws.Cells.Item(3, 6).Value = 1.744969725608826
$ws.Cells.Item(3, 7).Value = -0.8538760542869568
$ws.Cells.Item(3, 8).Value = -4.247409343719482

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "walkingToRunning"
$ws.Cells.Item(4, 3).Value = -3.206811956737363
$ws.Cells.Item(4, 4).Value = -14.22782378611358
$ws.Cells.Item(4, 5).Value = -3.585576596467396
$ws.Cells.Item(4, 6).Value = 4.650478363037109
$ws.Cells.Item(4, 7).Value = 11.04139995574951
$ws.Cells.Item(4, 8).Value = -3.084672927856445

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "walkingToRunning"
$ws.Cells.Item(5, 3).Value = -12.97183887854877
$ws.Cells.Item(5, 4).Value = -19.3549799711808
$ws.Cells.Item(5, 5).Value = -15.49820016778038
$ws.Cells.Item(5, 6).Value = -3.62076735496521
$ws.Cells.Item(5, 7).Value = -5.12558650970459
$ws.Cells.Item(5, 8).Value = 7.287806987762451

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "walkingToRunning"
$ws.Cells.Item(6, 3).Value = 6.874614751857277
$ws.Cells.Item(6, 4).Value = -17.72338112540867
$ws.Cells.Item(6, 5).Value = -7.79767358821371
$ws.Cells.Item(6, 6).Value = -5.951565742492676
$ws.Cells.Item(6, 7).Value = 1.577590107917786
$ws.Cells.Item(6, 8).Value = 6.344515323638916

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "walkingToRunning"
$ws.Cells.Item(7, 3).Value = 27.40819798345139
$ws.Cells.Item(7, 4).Value = -14.00006047539093
$ws.Cells.Item(7, 5).Value = 1.457801404206576
$ws.Cells.Item(7, 6).Value = -1.728857636451721
$ws.Cells.Item(7, 7).Value = 8.543354034423828
$ws.Cells.Item(7, 8).Value = 3.901864051818848

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "walkingToRunning"
$ws.Cells.Item(8, 3).Value = 40.11718940734848
$ws.Cells.Item(8, 4).Value = -4.808408617973335
$ws.Cells.Item(8, 5).Value = 16.45111835002897
$ws.Cells.Item(8, 6).Value = 1.824864625930786
$ws.Cells.Item(8, 7).Value = -0.6525410413742065
$ws.Cells.Item(8, 8).Value = -3.505985021591187

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "walkingToRunning"
$ws.Cells.Item(9, 3).Value = -6.829358805780815
$ws.Cells.Item(9, 4).Value = -18.71223431048175
$ws.Cells.Item(9, 5).Value = 7.999558770138284
$ws.Cells.Item(9, 6).Value = 5.772734642028809
$ws.Cells.Item(9, 7).Value = -3.56550669670105
$ws.Cells.Item(9, 8).Value = -3.791476011276245

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "walkingToRunning"
$ws.Cells.Item(10, 3).Value = -8.796926493230012
$ws.Cells.Item(10, 4).Value = -58.53746407446633
$ws.Cells.Item(10, 5).Value = 6.315518513969705
$ws.Cells.Item(10, 6).Value = 4.579105854034424
$ws.Cells.Item(10, 7).Value = 3.013566493988037
$ws.Cells.Item(10, 8).Value = -2.666023969650269

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "walkingToRunning"
$ws.Cells.Item(11, 3).Value = 3.485328860904898
$ws.Cells.Item(11, 4).Value = -5.610776631728454
$ws.Cells.Item(11, 5).Value = -5.665500184764032
$ws.Cells.Item(11, 6).Value = -5.41840124130249
$ws.Cells.Item(11, 7).Value = -4.322909832000732
$ws.Cells.Item(11, 8).Value = 2.317283153533936

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "walkingToRunning"
$ws.Cells.Item(12, 3).Value = -16.32176755822241
$ws.Cells.Item(12, 4).Value = -11.76148359671859
$ws.Cells.Item(12, 5).Value = -5.036116931749546
$ws.Cells.Item(12, 6).Value = -11.62250137329102
$ws.Cells.Item(12, 7).Value = -2.614225625991821
$ws.Cells.Item(12, 8).Value = 11.06896305084228

$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "walkingToRunning"
$ws.Cells.Item(13, 3).Value = 16.76886541947075
$ws.Cells.Item(13, 4).Value = -55.28290149439937
$ws.Cells.Item(13, 5).Value = 16.80483585855235
$ws.Cells.Item(13, 6).Value = 7.474494457244873
$ws.Cells.Item(13, 7).Value = -2.812897443771362
$ws.Cells.Item(13, 8).Value = -3.427688121795654

$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "walkingToRunning"
$ws.Cells.Item(14, 3).Value = -17.68672215420252
$ws.Cells.Item(14, 4).Value = -8.062301013780889
$ws.Cells.Item(14, 5).Value = 3.262691987597485
$ws.Cells.Item(14, 6).Value = 1.629921197891235
$ws.Cells.Item(14, 7).Value = -8.41618824005127
$ws.Cells.Item(14, 8).Value = 3.490139245986938

$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "walkingToRunning"
$ws.Cells.Item(15, 3).Value = 6.037598153819246
$ws.Cells.Item(15, 4).Value = -26.96122758284862
$ws.Cells.Item(15, 5).Value = 22.58520386530009
$ws.Cells.Item(15, 6).Value = 8.126436233520508
$ws.Cells.Item(15, 7).Value = -4.144477844238281
$ws.Cells.Item(15, 8).Value = -0.5264403820037842

$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "walkingToRunning"
$ws.Cells.Item(16, 3).Value = -27.64734548071142
$ws.Cells.Item(16, 4).Value = -30.63590854147206
$ws.Cells.Item(16, 5).Value = -13.72450681354674
$ws.Cells.Item(16, 6).Value = 1.081842660903931
$ws.Cells.Item(16, 7).Value = -4.525310039520264
$ws.Cells.Item(16, 8).Value = -0.875314474105835

$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "walkingToRunning"
$ws.Cells.Item(17, 3).Value = -0.3365890254143977
$ws.Cells.Item(17, 4).Value = -9.820995450019772
$ws.Cells.Item(17, 5).Value = -11.5744883495828
$ws.Cells.Item(17, 6).Value = -1.842840909957886
$ws.Cells.Item(17, 7).Value = -1.715674996376038
$ws.Cells.Item(17, 8).Value = -1.691839694976807

$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "walkingToRunning"
$ws.Cells.Item(18, 3).Value = -9.280692992003011
$ws.Cells.Item(18, 4).Value = -14.45616371476128
$ws.Cells.Item(18, 5).Value = -11.51578338249869
$ws.Cells.Item(18, 6).Value = -15.37276458740234
$ws.Cells.Item(18, 7).Value = -16.33429908752441
$ws.Cells.Item(18, 8).Value = 3.632352113723755

$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "walkingToRunning"
$ws.Cells.Item(19, 3).Value = -4.598872557930385
$ws.Cells.Item(19, 4).Value = -17.49396556356679
$ws.Cells.Item(19, 5).Value = 6.897694048674124
$ws.Cells.Item(19, 6).Value = 7.317900657653809
$ws.Cells.Item(19, 7).Value = -9.003682136535645
$ws.Cells.Item(19, 8).Value = -1.13044536113739

$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "walkingToRunning"
$ws.Cells.Item(20, 3).Value = -11.20024363890936
$ws.Cells.Item(20, 4).Value = 1.505196239637264
$ws.Cells.Item(20, 5).Value = 13.08335323955702
$ws.Cells.Item(20, 6).Value = 0.0980709120631218
$ws.Cells.Item(20, 7).Value = -6.593521118164063
$ws.Cells.Item(20, 8).Value = 0.4514724016189575

$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "walkingToRunning"
$ws.Cells.Item(21, 3).Value = -36.36955037324309
$ws.Cells.Item(21, 4).Value = -40.35419501428967
$ws.Cells.Item(21, 5).Value = 41.72007127430098
$ws.Cells.Item(21, 6).Value = 9.948039054870604
$ws.Cells.Item(21, 7).Value = 6.794189929962158
$ws.Cells.Item(21, 8).Value = 0.0589224398136138

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "walkingToRunning"
$ws.Cells.Item(22, 3).Value = -64.13763353098994
$ws.Cells.Item(22, 4).Value = -54.61329487095708
$ws.Cells.Item(22, 5).Value = 20.88871420984682
$ws.Cells.Item(22, 6).Value = -1.0385662317276
$ws.Cells.Item(22, 7).Value = -3.36523699760437
$ws.Cells.Item(22, 8).Value = 1.002214074134827

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "walkingToRunning"
$ws.Cells.Item(23, 3).Value = -48.54061980869488
$ws.Cells.Item(23, 4).Value = -15.36169694817591
$ws.Cells.Item(23, 5).Value = -18.69812476116687
$ws.Cells.Item(23, 6).Value = -2.924617052078247
$ws.Cells.Item(23, 7).Value = 0.8201870322227478
$ws.Cells.Item(23, 8).Value = 1.637111783027649

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "walkingToRunning"
$ws.Cells.Item(24, 3).Value = -9.763155895730725
$ws.Cells.Item(24, 4).Value = -17.02286973725198
$ws.Cells.Item(24, 5).Value = -7.877094351727034
$ws.Cells.Item(24, 6).Value = -9.886119842529297
$ws.Cells.Item(24, 7).Value = -17.29250526428223
$ws.Cells.Item(24, 8).Value = -0.4699813723564148

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "walkingToRunning"
$ws.Cells.Item(25, 3).Value = -12.50246284319005
$ws.Cells.Item(25, 4).Value = -25.3144741265671
$ws.Cells.Item(25, 5).Value = 2.428819822228345
$ws.Cells.Item(25, 6).Value = 5.15168571472168
$ws.Cells.Item(25, 7).Value = -2.607833862304688
$ws.Cells.Item(25, 8).Value = -3.956591844558716

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "walkingToRunning"
$ws.Cells.Item(26, 3).Value = -12.77425661294355
$ws.Cells.Item(26, 4).Value = -3.042295404102455
$ws.Cells.Item(26, 5).Value = 24.03419656857195
$ws.Cells.Item(26, 6).Value = -0.8223176002502441
$ws.Cells.Item(26, 7).Value = -5.908023357391357
$ws.Cells.Item(26, 8).Value = 0.7502790689468384

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "walkingToRunning"
$ws.Cells.Item(27, 3).Value = -11.1896470111351
$ws.Cells.Item(27, 4).Value = -11.74091952780091
$ws.Cells.Item(27, 5).Value = 15.447055526402
$ws.Cells.Item(27, 6).Value = 7.825499057769775
$ws.Cells.Item(27, 7).Value = -0.2695784866809845
$ws.Cells.Item(27, 8).Value = -1.960819005966186

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "walkingToRunning"
$ws.Cells.Item(28, 3).Value = -33.40386452882192
$ws.Cells.Item(28, 4).Value = -74.93744767230478
$ws.Cells.Item(28, 5).Value = 35.67094943834376
$ws.Cells.Item(28, 6).Value = 2.197174549102783
$ws.Cells.Item(28, 7).Value = -1.826995134353638
$ws.Cells.Item(28, 8).Value = -2.098770618438721

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "walkingToRunning"
$ws.Cells.Item(29, 3).Value = -55.71352932764142
$ws.Cells.Item(29, 4).Value = -17.14117759207022
$ws.Cells.Item(29, 5).Value = -22.77800974638568
$ws.Cells.Item(29, 6).Value = -1.679322838783264
$ws.Cells.Item(29, 7).Value = 2.35203742980957
$ws.Cells.Item(29, 8).Value = 2.066413402557373

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "walkingToRunning"
$ws.Cells.Item(30, 3).Value = -7.357292139011771
$ws.Cells.Item(30, 4).Value = -6.564726891724789
$ws.Cells.Item(30, 5).Value = -1.292621791362762
$ws.Cells.Item(30, 6).Value = -7.180881023406982
$ws.Cells.Item(30, 7).Value = -6.29631233215332
$ws.Cells.Item(30, 8).Value = 12.63969612121582

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "walkingToRunning"
$ws.Cells.Item(31, 3).Value = -26.65620994567871
$ws.Cells.Item(31, 4).Value = -36.59538269042969
$ws.Cells.Item(31, 5).Value = -1.42856240272522
$ws.Cells.Item(31, 6).Value = 0.2972753643989563
$ws.Cells.Item(31, 7).Value = 7.379552841186523
$ws.Cells.Item(31, 8).Value = -0.9200555682182312
